$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.074.17"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "'  +0.17%  "
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = "'1.833.27"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = "'  -0.01%  "
$ws.Range('E3').ClearFormats()
$ws.Range('D4').Value = "'0.9999"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = "'  +0.16%  "
$ws.Range('E4').ClearFormats()
$ws.Range('D5').Value = "'243.59"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "'  +0.57%  "
$ws.Range('E5').ClearFormats()
$ws.Range('D6').Value = "'0.6271"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "'  +0.07%  "
$ws.Range('E6').ClearFormats()
$ws.Range('D7').Value = "'1.001"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = "'  +0.19%  "
$ws.Range('E7').ClearFormats()
$ws.Range('D8').Value = "'0.07516"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = "'  -1.32%  "
$ws.Range('E8').ClearFormats()
$ws.Range('D9').Value = "'0.2922"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "'  -0.18%  "
$ws.Range('E9').ClearFormats()
$ws.Range('D10').Value = "'23.19"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = "'  +2.80%  "
$ws.Range('E10').ClearFormats()
$ws.Range('D11').Value = "'0.07676"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = "'  -0.48%  "
$ws.Range('E11').ClearFormats()
$ws.Range('D12').Value = "'1.832.74"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "'  -0.10%  "
$ws.Range('E12').ClearFormats()
$ws.Range('D13').Value = "'4.997"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "'  +0.87%  "
$ws.Range('E13').ClearFormats()
$ws.Range('D14').Value = "'0.6671"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = "'  +0.23%  "
$ws.Range('E14').ClearFormats()
$ws.Range('E15').Value = "'  -0.07%  "
$ws.Range('E15').ClearFormats()
$ws.Range('D16').Value = "'0.000009375"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = "'  -7.77%  "
$ws.Range('E16').ClearFormats()
$ws.Range('D17').Value = "'5.983"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = "'  -1.15%  "
$ws.Range('E17').ClearFormats()
$ws.Range('D18').Value = "'29.093.89"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "'  +0.20%  "
$ws.Range('E18').ClearFormats()
$ws.Range('D19').Value = "'2.075.62"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "'  -0.47%  "
$ws.Range('E19').ClearFormats()
$ws.Range('D20').Value = "'12.58"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = "'  +1.79%  "
$ws.Range('E20').ClearFormats()
$ws.Range('D21').Value = "'223.44"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "'  -1.38%  "
$ws.Range('E21').ClearFormats()
$ws.Range('D22').Value = "'1.003"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "'  +0.43%  "
$ws.Range('E22').ClearFormats()
$ws.Range('D23').Value = "'7.096"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = "'  -1.21%  "
$ws.Range('E23').ClearFormats()
$ws.Range('D24').Value = "'1.001"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "'  +0.15%  "
$ws.Range('E24').ClearFormats()
$ws.Range('D25').Value = "'159.75"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "'  +0.95%  "
$ws.Range('E25').ClearFormats()
$ws.Range('D26').Value = "'0.1396"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = "'  +1.47%  "
$ws.Range('E26').ClearFormats()
$ws.Range('D27').Value = "'8.488"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "'  -0.11%  "
$ws.Range('E27').ClearFormats()
$ws.Range('D28').Value = "'17.89"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "'  -0.06%  "
$ws.Range('E28').ClearFormats()
$ws.Range('D29').Value = "'1.496"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = "'  +0.39%  "
$ws.Range('E29').ClearFormats()
$ws.Range('D30').Value = "'0.05683"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = "'  +8.92%  "
$ws.Range('E30').ClearFormats()
$ws.Range('D31').Value = "'4.149"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = "'  +0.96%  "
$ws.Range('E31').ClearFormats()
$ws.Range('D32').Value = "'4.070"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = "'  +1.31%  "
$ws.Range('E32').ClearFormats()
$ws.Range('D33').Value = "'1.204"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = "'  +1.25%  "
$ws.Range('E33').ClearFormats()
$ws.Range('D34').Value = "'0.7430"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = "'  +0.96%  "
$ws.Range('E34').ClearFormats()
$ws.Range('D35').Value = "'1.839"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = "'  -0.29%  "
$ws.Range('E35').ClearFormats()
$ws.Range('D36').Value = "'1.139"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = "'  -0.08%  "
$ws.Range('E36').ClearFormats()
$ws.Range('D37').Value = "'2.669"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "'  -1.37%  "
$ws.Range('E37').ClearFormats()
$ws.Range('D38').Value = "'2.761"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = "'  +0.14%  "
$ws.Range('E38').ClearFormats()
$ws.Range('D39').Value = "'1.218.80"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "'  -1.84%  "
$ws.Range('E39').ClearFormats()
$ws.Range('D40').Value = "'0.01782"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "'  -0.17%  "
$ws.Range('E40').ClearFormats()
$ws.Range('D41').Value = "'6.515"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "'  +2.68%  "
$ws.Range('E41').ClearFormats()
$ws.Range('D42').Value = "'0.8928"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = "'  -0.41%  "
$ws.Range('E42').ClearFormats()
$ws.Range('E43').Value = "'  +0.27%  "
$ws.Range('E43').ClearFormats()
$ws.Range('D44').Value = "'101.92"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "'  +0.33%  "
$ws.Range('E44').ClearFormats()
$ws.Range('D45').Value = "'1.979.66"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "'  -0.15%  "
$ws.Range('E45').ClearFormats()
$ws.Range('E46').Value = "'  +2.15%  "
$ws.Range('E46').ClearFormats()
$ws.Range('D47').Value = "'65.69"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "'  +2.21%  "
$ws.Range('E47').ClearFormats()
$ws.Range('B48').Value = "'Mantle"
$ws.Range('B48').ClearFormats()
$ws.Range('C48').Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range('C48').ClearFormats()
$ws.Range('D48').Value = "'0.5085"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = "'  -0.36%  "
$ws.Range('E48').ClearFormats()
$ws.Range('B49').Value = "'XinFinNetwork"
$ws.Range('B49').ClearFormats()
$ws.Range('C49').Value = "'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range('C49').ClearFormats()
$ws.Range('D49').Value = "'0.07591"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = "'  +10.10%  "
$ws.Range('E49').ClearFormats()
$ws.Range('D50').Value = "'0.4074"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "'  +0.97%  "
$ws.Range('E50').ClearFormats()
$ws.Range('D51').Value = "'9.024"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = "'  +1.79%  "
$ws.Range('E51').ClearFormats()
